# Apply cryptos list update (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.376.22"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "3.506.99"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'591.41"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "'134.54"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +0.54%  "
$ws.Range("D9").Value = "'7.61"
$ws.Range("E9").Value = "  +6.03%  "
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").Value = "'0.389"
$ws.Range("E11").Value = "  +3.81%  "
$ws.Range("D12").Value = "4.105.19"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +1.02%  "
$ws.Range("E14").Value = "  +1.12%  "
$ws.Range("D15").Value = "3.510.18"
$ws.Range("E15").Value = "  +0.75%  "
$ws.Range("D16").Value = "'25.76"
$ws.Range("E16").Value = "  +3.53%  "
$ws.Range("D17").Value = "64.352.86"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "'391.82"
$ws.Range("E21").Value = "  +1.82%  "
$ws.Range("E22").Value = "  +2.92%  "
$ws.Range("D23").Value = "3.646.92"
$ws.Range("E23").Value = "  +0.73%  "
$ws.Range("D24").Value = "'74.48"
$ws.Range("E24").Value = "  +0.68%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("D26").Value = "'5.67"
$ws.Range("E26").Value = "  -0.95%  "
$ws.Range("E27").Value = "  +4.00%  "
$ws.Range("E28").Value = "  +0.39%  "
$ws.Range("E29").Value = "  +0.53%  "
$ws.Range("E30").Value = "  +2.37%  "
$ws.Range("D31").Value = "'8.22"
$ws.Range("E31").Value = "  +0.67%  "
$ws.Range("E32").Value = "  -4.50%  "
$ws.Range("E33").Value = "  +7.63%  "
$ws.Range("D34").Value = "3.533.99"
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "'23.38"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("E38").Value = "  +2.19%  "
$ws.Range("E39").Value = "  +2.12%  "
$ws.Range("D40").Value = "'166.77"
$ws.Range("E40").Value = "  +2.48%  "
$ws.Range("E41").Value = "  +1.55%  "
$ws.Range("E42").Value = "  +1.14%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.45"
$ws.Range("E44").Value = "  +1.77%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'25.00"
$ws.Range("E45").Value = "  -1.84%  "
$ws.Range("E46").Value = "  +0.30%  "
$ws.Range("E47").Value = "  -1.61%  "
$ws.Range("E48").Value = "  +1.07%  "
$ws.Range("D49").Value = "'0.912"
$ws.Range("E49").Value = "  +1.84%  "
$ws.Range("D50").Value = "2.374.89"
$ws.Range("E50").Value = "  -3.68%  "
$ws.Range("E51").Value = "  +0.65%  "
